# Applies the "simplificar as funcionalidades" edit:
#  - Sheet "Regras-de-Negocio": RN001-RN003 get new Nome/Descricao text
#    (premium/payment related rules replaced by account/goal related rules),
#    RN004's Nome/Descricao are cleared out.
#  - Sheet "Requisitos-Funcionais": RF002-RF008 get new text (premium related
#    requirements replaced by login/goal related requirements) and the old
#    RF009 entry is removed entirely (rows deleted).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: Regras-de-Negocio
# ---------------------------------------------------------------------------
$wsRN = $wb.Worksheets.Item("Regras-de-Negocio")

# RN001
$wsRN.Range("B4").Value = "Usuário;"
$wsRN.Range("B5").Value = "Cada úsuario terá uma conta atrelada a um perfil de condicionamentos iniciais;"

# RN002
$wsRN.Range("B8").Value = "Metas do usuário;"
$wsRN.Range("B9").Value = "O usuário deverá ter no mínimo uma meta atribuida para sí;"

# RN003
$wsRN.Range("B12").Value = "Acompanhamentos das metas;"
$wsRN.Range("B13").Value = "O usuário deve efetuar registros de acordo com a periodicidade de preenchimento estabelecida para cada meta por ele estabelecida;"

# RN004 (no longer used) - clear the Nome/Descricao cells
$wsRN.Range("B15").ClearContents()
$wsRN.Range("B16").ClearContents()
$wsRN.Range("B17").ClearContents()

# ---------------------------------------------------------------------------
# Sheet 2: Requisitos-Funcionais
# ---------------------------------------------------------------------------
$wsRF = $wb.Worksheets.Item("Requisitos-Funcionais")

# RF002 (was "Adquirir conta premium;")
$wsRF.Range("B7").Value = "Login e Logout no aplicativo;"

# RF003 (was "Validar pagamento da assinatura;")
$wsRF.Range("B10").Value = "Incluir metas para si;"

# RF004 (was "Cadastrar cartão de credito;")
$wsRF.Range("B13").Value = " Consultar suas metas;"

# RF005 (was "Incluir metas para si;")
$wsRF.Range("B16").Value = "Alterar suas metas;"

# RF006 (was " Consultar suas metas;")
$wsRF.Range("B19").Value = "Excluir metas para si;"

# RF007 (was "Alterar suas metas;")
$wsRF.Range("B22").Value = "Registrar dados das metas;"

# RF008 (was "Excluir metas para si;")
$wsRF.Range("B25").Value = "Acompanhar progressão de metas;"

# RF009 block is removed entirely: rows 26 (spacer), 27 (ID) and 28 (Nome)
$wsRF.Range("A26:B28").EntireRow.Delete() | Out-Null

# ---------------------------------------------------------------------------
# Final UI state: Regras-de-Negocio tab active with E7 selected,
# Requisitos-Funcionais left with B15 selected.
# ---------------------------------------------------------------------------
$wsRF.Range("B15").Select() | Out-Null
$wsRN.Activate() | Out-Null
$wsRN.Range("E7").Select() | Out-Null
